$d = $word.ActiveDocument

# Step 1: merge the two runs that were split around the old _GoBack bookmark
# into a single run, by replacing the whole paragraph
# ("注意：诊断后的模型...") with equivalent XML (bookmark removed, text merged).
$full37 = $d.Paragraphs.Item(37).Range
$xmlPara37 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Times" w:eastAsia="Times New Roman" w:hAnsi="Times" w:cs="Times New Roman"/>
      <w:kern w:val="0"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="eastAsia"/>
      <w:color w:val="FF6600"/>
    </w:rPr>
    <w:t>注意</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="eastAsia"/>
    </w:rPr>
    <w:t>：</w:t>
  </w:r>
  <w:r>
    <w:t>诊断后的模型需要进行调优，调优后的新模型需要重新进行诊断，这是一个反复迭代不断逼近的过程，需要不断地尝试，</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>进而达到最优状态</w:t>
  </w:r>
</w:p>
'@
$full37.InsertXML($xmlPara37)

# Step 2: replace the whole tail region (from the "朴素贝叶斯" heading through
# the end of the document) with the updated content: pPr-simplified
# paragraphs, the new "自然常数e" heading + table, and the trailing empty
# paragraphs.
$tailStart = $d.Paragraphs.Item(42).Range.Start
$tailRange = $d.Range($tailStart, $d.Content.End)
$xmlRegion = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="eastAsia"/>
      <w:color w:val="0000FF"/>
    </w:rPr>
    <w:t>朴素贝叶斯</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="eastAsia"/>
    </w:rPr>
    <w:t>：</w:t>
  </w:r>
</w:p>
<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:tblPr>
    <w:tblStyle w:val="a3"/>
    <w:tblW w:w="0" w:type="auto"/>
    <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
  </w:tblPr>
  <w:tblGrid>
    <w:gridCol w:w="8516"/>
  </w:tblGrid>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="8516" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
            <w:b/>
            <w:color w:val="FF6600"/>
          </w:rPr>
          <w:t>为什么非常实用</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t>：</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t>首先，它是一个逆概问题，揭示了</w:t>
        </w:r>
        <w:r>
          <w:t>P(X|Y)</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t>和</w:t>
        </w:r>
        <w:r>
          <w:t>P(Y|X)</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t>的相反方向的条件概率的转换问题。</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t>从贝叶斯公式的发现历史来看，其就是为了处理所谓“逆概”问题而诞生的。比如</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t xml:space="preserve">P(Y|X) </w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t>不能通过直接观测来得到结果，而</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t xml:space="preserve">P(X|Y) </w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t>却容易通过直接观测得到结果，就可以通过贝叶斯公式从间接地观测对象去推断不可直接观测的对象的情况。</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t>例如：</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t>引申一步，基于样本特征去判断其所属标签的概率不好求，但是基于已经搜集好的打上标签的样本（有监督），却可以直接统计属于同一标签的样本内部各个特征的概率分布。因此贝叶斯方法的理论视角适用于一切分类问题的求解。</w:t>
        </w:r>
      </w:p>
      <w:p/>
      <w:p/>
    </w:tc>
  </w:tr>
</w:tbl>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="eastAsia"/>
      <w:color w:val="0000FF"/>
    </w:rPr>
    <w:t>自然常数</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:color w:val="0000FF"/>
    </w:rPr>
    <w:t>e</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="eastAsia"/>
    </w:rPr>
    <w:t>：</w:t>
  </w:r>
</w:p>
<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:tblPr>
    <w:tblStyle w:val="a3"/>
    <w:tblW w:w="0" w:type="auto"/>
    <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
  </w:tblPr>
  <w:tblGrid>
    <w:gridCol w:w="8516"/>
  </w:tblGrid>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="8516" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:proofErr w:type="gramStart"/>
        <w:r>
          <w:t>e</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:vertAlign w:val="superscript"/>
          </w:rPr>
          <w:t>x</w:t>
        </w:r>
        <w:proofErr w:type="gramEnd"/>
        <w:r>
          <w:t xml:space="preserve"> =1+ x+ x</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:vertAlign w:val="superscript"/>
          </w:rPr>
          <w:t>2</w:t>
        </w:r>
        <w:r>
          <w:t xml:space="preserve">/2! + </w:t>
        </w:r>
        <w:proofErr w:type="gramStart"/>
        <w:r>
          <w:t>x</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:vertAlign w:val="superscript"/>
          </w:rPr>
          <w:t>3</w:t>
        </w:r>
        <w:proofErr w:type="gramEnd"/>
        <w:r>
          <w:t>/3! + …</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t>这个式子可以得到常数</w:t>
        </w:r>
        <w:r>
          <w:t>e</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t>的范围，可以得到</w:t>
        </w:r>
        <w:proofErr w:type="spellStart"/>
        <w:r>
          <w:t>possion</w:t>
        </w:r>
        <w:proofErr w:type="spellEnd"/>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t>分布的</w:t>
        </w:r>
        <w:r>
          <w:t>pdf.</w:t>
        </w:r>
        <w:bookmarkStart w:id="0" w:name="_GoBack"/>
        <w:bookmarkEnd w:id="0"/>
      </w:p>
    </w:tc>
  </w:tr>
</w:tbl>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@
$tailRange.InsertXML($xmlRegion)

Write-Host "done"
